$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = "AAF"
$ws.Range("A44").Value = "AAF"
